# Added analysis UK Pillar 2 SGTF data
#
# This script nudges the position/size of a handful of scatter-point
# ellipses (and the two size-legend bubbles + their labels) on the one
# slide of the deck. All of the affected shapes live inside the single
# top-level group shape that holds the whole plot.
#
# PowerPoint's Shape.Left/Top/Width/Height are expressed in points while
# the underlying OOXML stores EMUs (1 pt = 12700 EMU). The COM layer here
# rounds points -> EMU by truncation, so a literal "target_emu / 12700"
# can land one EMU short. We add a half-EMU (plus a tiny guard) worth of
# points before converting so the stored EMU value matches the target
# exactly.

function Get-ShapeById($collection, $id) {
    for ($i = 1; $i -le $collection.Count; $i++) {
        $item = $collection.Item($i)
        if ($item.Id -eq $id) {
            return $item
        }
    }
    return $null
}

$EmuPerPoint = 12700
$HalfEmuInPoints = (0.5 / $EmuPerPoint) + 0.000000001

function Emu-ToPoints($emu) {
    return ($emu / $EmuPerPoint) + $HalfEmuInPoints
}

function Set-ShapeRectEmu($shape, $offX, $offY, $extCx, $extCy) {
    $shape.Left = Emu-ToPoints $offX
    $shape.Top = Emu-ToPoints $offY
    $shape.Width = Emu-ToPoints $extCx
    $shape.Height = Emu-ToPoints $extCy
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$grp = $s.Shapes.Item(1)
$items = $grp.GroupItems

# pt14
Set-ShapeRectEmu (Get-ShapeById $items 14) 2248371 6070370 67453 67453

# pt15
Set-ShapeRectEmu (Get-ShapeById $items 15) 2830720 6075092 58010 58010

# pt16
Set-ShapeRectEmu (Get-ShapeById $items 16) 3403626 5686041 67453 67453

# pt17
Set-ShapeRectEmu (Get-ShapeById $items 17) 3981254 3380064 67453 67453

# pt18
Set-ShapeRectEmu (Get-ShapeById $items 18) 4546081 3254776 93055 93055

# pt19
Set-ShapeRectEmu (Get-ShapeById $items 19) 5119712 3198553 101049 101049

# pt20
Set-ShapeRectEmu (Get-ShapeById $items 20) 5692077 1661330 111573 111573

# pt47
Set-ShapeRectEmu (Get-ShapeById $items 47) 7629570 3578317 9271 9271

# pt49
Set-ShapeRectEmu (Get-ShapeById $items 49) 7603069 3771272 62273 62273

# tx52 (legend bubble label "10" -> "1", also reflows/shrinks its box)
$tx52 = Get-ShapeById $items 52
Set-ShapeRectEmu $tx52 7819849 3538959 67806 87630
$tx52.TextFrame.TextRange.Text = "1"

# tx53 (legend bubble label "50" -> "10", position/size unchanged)
$tx53 = Get-ShapeById $items 53
$tx53.TextFrame.TextRange.Text = "10"
